$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.252.03"
$ws.Range("E2").Value = "  -2.25%  "
$ws.Range("D3").Value = "1.878.35"
$ws.Range("E3").Value = "  -2.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.85"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4836"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2871"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06581"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.00%  "
$ws.Range("D10").Value = "1.884.32"
$ws.Range("E10").Value = "  -1.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.70"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07319"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.131"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.03"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6530"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.08%  "
$ws.Range("D16").Value = "30.223.67"
$ws.Range("E16").Value = "  -2.19%  "
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007731"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").Value = "2.132.50"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.364"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "193.53"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -7.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.099"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.251"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.12"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.98"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.908"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.432"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.256"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09075"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.001"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05056"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7115"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.094"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.698"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01777"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.633"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9210"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.036"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "105.77"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4259"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -6.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.767"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.380"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1311"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -7.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "64.93"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.895"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05748"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.61"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3811"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -7.12%  "
